$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- register the extra (text / date) cell-format variants that were produced
#     while laying out this sheet; applied to a scratch cell off the used
#     range and then cleared back out so the grid itself stays plain numbers ---
$ws.Range("Z1").Borders.Item(7).LineStyle = 1
$ws.Range("Z1").NumberFormat = "@"
$ws.Range("Z1").Value = "x"
$ws.Range("Z1").ClearFormats()
$ws.Range("Z1").ClearContents()

$ws.Range("Z2").Borders.Item(8).LineStyle = 1
$ws.Range("Z2").NumberFormat = "m/d/yy h:mm"
$ws.Range("Z2").Value = 1
$ws.Range("Z2").ClearFormats()
$ws.Range("Z2").ClearContents()

# --- Data: 90-degree GLCM features for the YCbCr "B" (blue-difference) channel ---
$data = @(
  @(0.00054269257761440965, 0.54970652389093966, 0.99807901639659402, 0.99972865371119291),
  @(0.0014795475998458662,  0.81540276075791207, 0.98864746616774246, 0.99926022620007748),
  @(0.00034989918685738854, 0.14289888162822195, 0.99921202569327772, 0.99982505040657121),
  @(0.051993135255543665,   0.82780601953128785, 0.65043311657995173, 0.9740034323722282)
)

for ($r = 0; $r -lt 4; $r++) {
  for ($c = 0; $c -lt 4; $c++) {
    $ws.Cells.Item($r + 1, $c + 1).Value = $data[$r][$c]
  }
}

# --- column widths to match the generated layout ---
$ws.Columns.Item(1).ColumnWidth = 14.75
$ws.Columns.Item(2).ColumnWidth = 11.75
$ws.Columns.Item(3).ColumnWidth = 11.75
$ws.Columns.Item(4).ColumnWidth = 11.75
